$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column I ("done by"): who reviewed / filled in each participant row ---
# Values are written in first-occurrence order so the shared-string table
# ends up built in the same order as the target workbook.
$ws.Range("I1").Value = "done by"
$ws.Range("I7").Value = "annebelle"
$ws.Range("I2").Value = "robin"
$ws.Range("I26").Value = "ellora"
$ws.Range("I20").Value = "tomas"
$ws.Range("I25").Value = "ellora?"

# --- Column J ("json?"): yes/no flag ---
$ws.Range("J1").Value = "json?"
$ws.Range("J25").Value = "no"
$ws.Range("J2").Value = "yes"

# --- Column K ("notes") ---
$ws.Range("K1").Value = "notes"
$ws.Range("K5").Value = "drop older version of 4"

# --- Remaining column I values ---
$ws.Range("I3").Value = "robin"
$ws.Range("I4").Value = "robin"
$ws.Range("I5").Value = "robin"
$ws.Range("I6").Value = "robin"
$ws.Range("I8").Value = "annebelle"
$ws.Range("I9").Value = "annebelle"
$ws.Range("I10").Value = "annebelle"
$ws.Range("I11").Value = "annebelle"
$ws.Range("I12").Value = "annebelle"
$ws.Range("I13").Value = "annebelle"
$ws.Range("I14").Value = "annebelle"
$ws.Range("I15").Value = "annebelle"
$ws.Range("I21").Value = "tomas"
$ws.Range("I22").Value = "tomas"
$ws.Range("I24").Value = "annebelle"
$ws.Range("I27").Value = "ellora?"
$ws.Range("I28").Value = "ellora?"

# --- Remaining column J values ---
$ws.Range("J3").Value = "yes"
$ws.Range("J4").Value = "yes"
$ws.Range("J5").Value = "yes"
$ws.Range("J6").Value = "yes"
$ws.Range("J7").Value = "yes"
$ws.Range("J8").Value = "yes"
$ws.Range("J9").Value = "yes"
$ws.Range("J10").Value = "yes"
$ws.Range("J11").Value = "yes"
$ws.Range("J12").Value = "yes"
$ws.Range("J13").Value = "yes"
$ws.Range("J14").Value = "yes"
$ws.Range("J15").Value = "yes"
$ws.Range("J20").Value = "yes"
$ws.Range("J21").Value = "yes"
$ws.Range("J22").Value = "yes"
$ws.Range("J24").Value = "yes"
$ws.Range("J26").Value = "yes"
$ws.Range("J27").Value = "no"
$ws.Range("J28").Value = "no"

# --- Column widths for the two newly introduced columns ---
$ws.Columns.Item(9).ColumnWidth = 12.666666666666666
$ws.Columns.Item(11).ColumnWidth = 21.5

# --- Normalize row heights for the data rows below the frozen header pane ---
for ($r = 16; $r -le 33; $r++) {
    $ws.Rows.Item($r).RowHeight = 20
}

# --- Leave the selection where the author last left it ---
$ws.Range("I24").Select() | Out-Null
